$d = $word.ActiveDocument

# --- Paragraph A: "... issue is attack entities ... destroyed when the laser collides ..."
# 1) "attack" -> "projectile"
$r = $d.Content.Duplicate
$r.Find.Execute("issue is attack entities", $true, $false, $false, $false, $false, $true, 1, $false, "issue is projectile entities", 1)

# 2) "laser" -> "projectile" (in "destroyed when the laser collides")
$r = $d.Content.Duplicate
$r.Find.Execute("destroyed when the laser collides", $true, $false, $false, $false, $false, $true, 1, $false, "destroyed when the projectile collides", 1)

# --- Paragraph B: "... before the level starts ..." and "The only overhead that is required is to ..."
# 3) "level" -> "gameplay"
$r = $d.Content.Duplicate
$r.Find.Execute("before the level starts", $true, $false, $false, $false, $false, $true, 1, $false, "before the gameplay starts", 1)

# 4) "The only overhead that is required is to" -> "Instead, the only overhead is to"
$r = $d.Content.Duplicate
$r.Find.Execute("The only overhead that is required is to activate/deactivate", $true, $false, $false, $false, $false, $true, 1, $false, "Instead, the only overhead is to activate/deactivate", 1)

# --- Paragraph C: "Instead of creating an Object Pool ... dragging them all into the scene, there is a singular" / "job is to give out references to the appropriate"
# 5) "scene" -> "Scene" (capitalize, in "dragging them all into the scene, there is a singular")
$r = $d.Content.Duplicate
$r.Find.Execute("dragging them all into the scene, there is a singular", $true, $false, $false, $false, $false, $true, 1, $false, "dragging them all into the Scene, there is a singular", 1)

# 6) "job is to give out references to the appropriate" -> "job is to give out references of the appropriate"
$r = $d.Content.Duplicate
$r.Find.Execute("job is to give out references to the appropriate", $true, $false, $false, $false, $false, $true, 1, $false, "job is to give out references of the appropriate", 1)

# --- Paragraph D: "Each Object Pool Prefab is then added to the singular Object Pool Manager Prefab.  The Object Pool Manager Prefab, is what is dragged into the scene."
# 7) Combined rewrite
$r = $d.Content.Duplicate
$r.Find.Execute("Each Object Pool Prefab is then added to the singular Object Pool Manager Prefab.  The Object Pool Manager Prefab, is what is dragged into the scene.", $true, $false, $false, $false, $false, $true, 1, $false, "Each Object Pool Prefab is added to the Object Pool Manager Prefab.  A singular instance of the Object Pool Manager Prefab is what is dragged into the scene.", 1)
